# Commit: "changed date on apply doc to 15 March"
# The application deadline date is moved from 20 March 2020 to 15 March 2020
# within the bold "To apply please..." sentence. Formatting (bold,
# Helvetica Neue, size 24) is preserved because Find/Replace only rewrites
# the matched text inside the existing run(s).

$d = $word.ActiveDocument

$found = $d.Content.Find.Execute(
    "by midnight 20 March 2020", $true, $false, $false, $false, $false,
    $true, 1, $false, "by midnight 15 March 2020", 2)

if (-not $found) {
    throw "Could not find the '20 March 2020' deadline text to update."
}
